# Weekly update: insert a new price record at the top of the data block
# (row 28), pushing the existing rows 28-84 down to 29-85.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 28; Excel shifts rows 28:84 down to 29:85
# and carries the row formatting (e.g. the date-style on column D) along.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new observation.
$ws.Range("A28").Value = 7
$ws.Range("B28").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C28").Value = "Ñuble"
$ws.Range("D28").Value = 45082
$ws.Range("E28").Value = 16
$ws.Range("F28").Value = 100112013
$ws.Range("G28").Value = "Alcachofa"
$ws.Range("H28").Value = "Argentina(o)"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 30
$ws.Range("K28").Value = 17000
$ws.Range("L28").Value = 17000
$ws.Range("M28").Value = 17000
$ws.Range("N28").Value = "`$/caja 50 unidades"
$ws.Range("O28").Value = "Provincia de Limarí"
$ws.Range("P28").Value = 340
$ws.Range("Q28").Value = 50
$ws.Range("R28").Value = "Hortaliza"
